# Replace the old 5-column "Source/Target mapping" table with a new,
# single-column list of validation rules (plus leftover time-formatted
# cells in column C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous table's values (A1:E6) while leaving any cell
# formatting (e.g. the time-format style on column C) intact.
$ws.Range("A1:E6").ClearContents()

# New validation-rule descriptions in column A.
$ws.Range("A1").Value = "Validate if the comments are provided, If present Print YES else Print NO"
$ws.Range("A2").Value = "Validate if the Snowflake SQL Keywords are in Upper Case, If present Print YES else Print NO"
$ws.Range("A3").Value = "Validate if there is a date hardcoded, If present Print YES else Print NO"

# Column A now holds much longer text, so it needs to be widened.
$ws.Columns.Item(1).ColumnWidth = 46.498697916666664

# Move the active cell/selection to reflect the new, shorter table.
$ws.Range("B7").Select()
